# Applies the cryptos-list price/volume refresh described in the commit
# "Updated cryptos list on Thu Apr 11 09:32:39 UTC 2024 with GitHub Actions".
# Column D values that look like plain numbers get a leading apostrophe so Excel
# stores them as text (matching the original text-formatted price cells) instead
# of silently converting them into floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.843.50'
$ws.Range("E2").Value = '  +2.56%  '
$ws.Range("D3").Value = '3.591.08'
$ws.Range("E3").Value = '  +2.00%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '''599.14'
$ws.Range("E5").Value = '  +1.35%  '
$ws.Range("D6").Value = '''173.61'
$ws.Range("E6").Value = '  +1.08%  '
$ws.Range("D7").Value = '3.583.93'
$ws.Range("E7").Value = '  +1.93%  '
$ws.Range("E8").Value = '  +0.66%  '
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("D10").Value = '''0.200'
$ws.Range("E10").Value = '  +5.69%  '
$ws.Range("D11").Value = '''7.49'
$ws.Range("E11").Value = '  +7.77%  '
$ws.Range("D12").Value = '''0.591'
$ws.Range("E12").Value = '  +1.42%  '
$ws.Range("D13").Value = '''46.91'
$ws.Range("D14").Value = '''0.0000279'
$ws.Range("E14").Value = '  +0.74%  '
$ws.Range("D15").Value = '4.170.99'
$ws.Range("E15").Value = '  +2.02%  '
$ws.Range("D16").Value = '''8.44'
$ws.Range("E16").Value = '  -0.56%  '
$ws.Range("D17").Value = '''614.42'
$ws.Range("E17").Value = '  -1.36%  '
$ws.Range("D18").Value = '3.580.45'
$ws.Range("E18").Value = '  +1.62%  '
$ws.Range("D19").Value = '70.887.30'
$ws.Range("E19").Value = '  +2.45%  '
$ws.Range("E20").Value = '  -1.06%  '
$ws.Range("D21").Value = '''17.50'
$ws.Range("E21").Value = '  +0.33%  '
$ws.Range("D22").Value = '''0.889'
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("D23").Value = '''9.28'
$ws.Range("E23").Value = '  -16.86%  '
$ws.Range("D24").Value = '''16.00'
$ws.Range("E24").Value = '  +0.48%  '
$ws.Range("D25").Value = '''97.32'
$ws.Range("E25").Value = '  +0.27%  '
$ws.Range("D26").Value = '''3.78'
$ws.Range("E26").Value = '  -1.71%  '
$ws.Range("D28").Value = '''2.66'
$ws.Range("E28").Value = '  +0.66%  '
$ws.Range("D29").Value = '''34.02'
$ws.Range("E29").Value = '  +3.73%  '
$ws.Range("D30").Value = '''9.22'
$ws.Range("E30").Value = '  -0.71%  '
$ws.Range("D31").Value = '''8.47'
$ws.Range("E31").Value = '  -0.84%  '
$ws.Range("D32").Value = '''3.08'
$ws.Range("E32").Value = '  -2.00%  '
$ws.Range("D33").Value = '''7.25'
$ws.Range("E33").Value = '  +4.43%  '
$ws.Range("E34").Value = '  -1.54%  '
$ws.Range("D35").Value = '''643.73'
$ws.Range("E35").Value = '  +0.94%  '
$ws.Range("D36").Value = '''3.72'
$ws.Range("E36").Value = '  +6.36%  '
$ws.Range("E37").Value = '  -1.24%  '
$ws.Range("D38").Value = '''10.87'
$ws.Range("E38").Value = '  +0.73%  '
$ws.Range("D39").Value = '''0.0483'
$ws.Range("E39").Value = '  +5.86%  '
$ws.Range("D40").Value = '''57.33'
$ws.Range("E40").Value = '  -0.10%  '
$ws.Range("E41").Value = '  +0.02%  '
$ws.Range("D42").Value = '''0.142'
$ws.Range("E42").Value = '  +4.85%  '
$ws.Range("D43").Value = '3.411.06'
$ws.Range("E43").Value = '  +0.79%  '
$ws.Range("D44").Value = '''0.325'
$ws.Range("E44").Value = '  -0.90%  '
$ws.Range("D45").Value = '0.0₃0720'
$ws.Range("E45").Value = '  +2.80%  '
$ws.Range("D46").Value = '''33.03'
$ws.Range("E46").Value = '  +0.45%  '
$ws.Range("D47").Value = '''2.97'
$ws.Range("E47").Value = '  +7.48%  '
$ws.Range("D48").Value = '''2.67'
$ws.Range("E48").Value = '  +4.97%  '
$ws.Range("E49").Value = '  +0.54%  '
$ws.Range("D50").Value = '''132.90'
$ws.Range("E50").Value = '  -0.30%  '
$ws.Range("E51").Value = '  -0.09%  '
